$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the value for the existing row 17 (H17)
$ws.Range("H17").Value = 92.875

# Insert 4 new rows at position 19 (this pushes the current rows 19-31,
# the "lif" block and everything after, down to rows 23-35).
$ws.Rows("19:22").Insert()

# Row 18 gains a new Epochs value and updated accuracy
$ws.Range("G18").Value = 100
$ws.Range("H18").Value = 93.875

# Fill the 4 newly inserted rows (19-22) with the new "Valve/adam" results
# that use 100 epochs instead of 200. Columns A and B stay blank, matching
# the existing "blank" cell pattern used elsewhere in the sheet (copy an
# already-blank cell so the blank cell itself is materialized, same as
# the rest of the table, instead of simply leaving it unset).
$ws.Range("A4").Copy($ws.Range("A19"))
$ws.Range("A4").Copy($ws.Range("B19"))
$ws.Range("C19").Value = "Valve"
$ws.Range("D19").Value = "adam"
$ws.Range("E19").Value = 0.0005
$ws.Range("F19").Value = 256
$ws.Range("G19").Value = 100
$ws.Range("H19").Value = 93.9375

$ws.Range("A4").Copy($ws.Range("A20"))
$ws.Range("A4").Copy($ws.Range("B20"))
$ws.Range("C20").Value = "Valve"
$ws.Range("D20").Value = "adam"
$ws.Range("E20").Value = 0.0005
$ws.Range("F20").Value = 256
$ws.Range("G20").Value = 100
$ws.Range("H20").Value = 93.3125

$ws.Range("A4").Copy($ws.Range("A21"))
$ws.Range("A4").Copy($ws.Range("B21"))
$ws.Range("C21").Value = "Valve"
$ws.Range("D21").Value = "adam"
$ws.Range("E21").Value = 0.0005
$ws.Range("F21").Value = 256
$ws.Range("G21").Value = 100
$ws.Range("H21").Value = 94.3125

$ws.Range("A4").Copy($ws.Range("A22"))
$ws.Range("A4").Copy($ws.Range("B22"))
$ws.Range("C22").Value = "Valve"
$ws.Range("D22").Value = "adam"
$ws.Range("E22").Value = 0.0005
$ws.Range("F22").Value = 256
$ws.Range("G22").Value = 100
$ws.Range("H22").Value = 92.375

Write-Host "Edit complete"
